# Apply updated dSF (column F) values to the walker_taijuan save-data sheet.
# These represent re-pulled data for the "push all data, mean calculation" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 7
    3  = -2
    4  = -1
    5  = 3
    6  = 6
    7  = -4
    9  = -1
    11 = 5
    12 = -6
    13 = 5
    14 = 2
    17 = 3
    19 = 3
    20 = 3
    21 = 1
    22 = -4
    23 = -1
    24 = 2
    25 = 2
    26 = -1
    27 = 3
    28 = -7
    31 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
